$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "Periodo Mora" row (2507) for CARLOS ALFREDO REYES MARTINEZ ---
# The existing block of rows for this worker (rows 18-22: periods 2506,2505,2504,2503,2502)
# gains a new most-recent period "2507" on top; everything below (including the later
# blank rows and the signature block at rows 27-28) shifts down by one row.
$ws.Rows("18").Insert()

# Copy formatting only from the row below (the old row18, now row19) so the new row
# matches the table's interior row style (border/fill/font) instead of default formatting.
$ws.Range("B19:J19").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)

# Fill in the new row's values
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1052990904"
$ws.Range("D18").Value = "CARLOS ALFREDO REYES MARTINEZ"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 52000
$ws.Range("G18").Value = 1300000

# --- Update the summary figures above the table ---
# VALOR MORA total
$ws.Range("E11").Value = 416000

# Cant. Periodos (count of periods per worker) increased because of the new period
$ws.Range("F13").Value = 6
